# Applies the SectorGroup.xlsx fix: for every row, the values in columns
# E (group-code), F (group-name) and G (category-name) were shifted by
# one position (E->F, F->G, G->E), i.e. after the edit:
#   new E = old G   (category-name)
#   new F = old E   (group-code)
#   new G = old F   (group-name)
# This applies uniformly to the header row as well as every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value()
    $fVal = $ws.Cells.Item($r, 6).Value()
    $gVal = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 5).Value = $gVal
    $ws.Cells.Item($r, 6).Value = $eVal
    $ws.Cells.Item($r, 7).Value = $fVal
}
